$d = $word.ActiveDocument

# --- 1) Merge the "Resumen" paragraph runs (remove the spell-check run
#        splits around "anilisis" / "triqui") into one single run. ---
$d.Content.Find.Execute(
    "El presente documento contiene el anilisis y diseño de un triqui distribuido bajo arquitectura P2P.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El presente documento contiene el anilisis y diseño de un triqui distribuido bajo arquitectura P2P.",
    2) | Out-Null

# --- 2) Merge the "Palabras clave" paragraph runs (remove the
#        spell-check run split around "triqui") into one single run. ---
$d.Content.Find.Execute(
    "Sistema distribuido, P2P, triqui, java, RMI.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sistema distribuido, P2P, triqui, java, RMI.",
    2) | Out-Null

# --- 3) Add the two new bibliography hyperlinks after the existing one,
#        plus one extra trailing blank paragraph. ---
#
# NOTE: Hyperlinks.Add always inserts its new run *before* whatever is
# already in the target range rather than replacing it, and creating a
# brand-new trailing paragraph at the very end of the story leaves a
# stray empty run behind. To keep the result byte-for-byte clean we grow
# the tail using real placeholder text (so any paragraph split happens
# between two non-empty runs) and then delete the placeholder text once
# it is no longer needed.

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertAfter("PH1")

$p = $d.Paragraphs.Item($count)
$h1 = $d.Hyperlinks.Add($p.Range, "http://cs.mty.itesm.mx/profesores/raul.perez/DAD/ejercicios/RMI/RMInetbeans.pdf", "", "", "http://cs.mty.itesm.mx/profesores/raul.perez/DAD/ejercicios/RMI/RMInetbeans.pdf")
$h1.Range.Style = "Hipervnculo"
$p = $d.Paragraphs.Item($count)
$delStart = $h1.Range.End
$delEnd = $p.Range.End - 1
if ($delEnd -gt $delStart) {
    $d.Range($delStart, $delEnd).Delete()
}

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertAfter("`rPH2")

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$h2 = $d.Hyperlinks.Add($p.Range, "http://download.java.net/jdk8/docs/technotes/guides/rmi/hello/hello-world.html", "", "", "http://download.java.net/jdk8/docs/technotes/guides/rmi/hello/hello-world.html")
$h2.Range.Style = "Hipervnculo"
$p = $d.Paragraphs.Item($count)
$delStart = $h2.Range.End
$delEnd = $p.Range.End - 1
if ($delEnd -gt $delStart) {
    $d.Range($delStart, $delEnd).Delete()
}

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertAfter("`rPH3")

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$delStart = $p.Range.Start
$delEnd = $p.Range.End - 1
if ($delEnd -gt $delStart) {
    $d.Range($delStart, $delEnd).Delete()
}
